$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 247 previously only carried A/D/E (the "pending" last row of the
# series). Now that a newer day exists, it gets its B/C figures filled in.
$ws.Cells.Item(247, 2).Value = 187
$ws.Cells.Item(247, 3).Value = 628

# Continue the daily date series with four more business days.
$newDates = @("04-09-2021", "05-09-2021", "06-09-2021", "07-09-2021")
$startRow = 248
$endRow = $startRow + $newDates.Length - 1

# Force column A to be read as text for the new rows so the dd-mm-yyyy
# strings are stored verbatim (as shared strings) instead of being
# auto-converted into date serial numbers, then drop the explicit
# number-format style again so the cells end up unstyled, same as the
# rest of the column.
$dateRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 1))
$dateRange.NumberFormat = "@"
for ($i = 0; $i -lt $newDates.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $newDates[$i]
}
$dateRange.Style = "Normal"

# D/E are populated for every new row.
for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Cells.Item($r, 4).Value = 3940
    $ws.Cells.Item($r, 5).Value = 30
}

# B/C are populated for every new row except the newest (last) one, which
# mirrors the "pending" pattern row 247 had before this update.
for ($r = $startRow; $r -lt $endRow; $r++) {
    $ws.Cells.Item($r, 2).Value = 187
    $ws.Cells.Item($r, 3).Value = 628
}
